$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45931, 2.245, 0.699),
    @(45931.01041666666, 8.002000000000001, 0.124),
    @(45931.02083333334, 1.875, 10.041),
    @(45931.03125, 5.481, 1.593),
    @(45931.04166666666, 0.08500000000000001, 13.576),
    @(45931.05208333334, 0, 17.309),
    @(45931.0625, 0, 14.541),
    @(45931.07291666666, 0, 21.32),
    @(45931.08333333334, 0, 22.491),
    @(45931.09375, 0, 11.023),
    @(45931.10416666666, 0, 15.427),
    @(45931.11458333334, 0, 16.463),
    @(45931.125, 0, 22.727),
    @(45931.13541666666, 0, 13.189),
    @(45931.14583333334, 0, 23.344),
    @(45931.15625, 0, 19.185),
    @(45931.16666666666, 2.584, 0.458),
    @(45931.17708333334, 3.65, 0.26),
    @(45931.1875, 0.297, 0.604),
    @(45931.19791666666, 1.437, 2.276),
    @(45931.20833333334, 0, 26.605),
    @(45931.21875, 0, 30.381),
    @(45931.22916666666, 0.028, 17.784),
    @(45931.23958333334, 0.211, 5.765),
    @(45931.25, 0.042, 17.331),
    @(45931.26041666666, 0.06, 9.65),
    @(45931.27083333334, 0.149, 5.902),
    @(45931.28125, 0.18, 5.68),
    @(45931.29166666666, 0.443, 2.188),
    @(45931.30208333334, 0.255, 5.102),
    @(45931.3125, 1.489, 7.201),
    @(45931.32291666666, 0.004, 16.089),
    @(45931.33333333334, 6.073, 8.074),
    @(45931.34375, 0, 18.616),
    @(45931.35416666666, 0.13, 2.751),
    @(45931.36458333334, 0, 1.504),
    @(45931.375, 2.861, 1.626),
    @(45931.38541666666, 0, 6.351),
    @(45931.39583333334, 0.076, 2.401),
    @(45931.40625, 0, 14.783),
    @(45931.41666666666, 0, 17.8),
    @(45931.42708333334, 0, 17.146)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

# Ensure the date/time number format is applied to all timestamp cells,
# including newly added rows beyond the original A1:C29 range.
$ws.Range("A2:A43").NumberFormat = "YYYY-MM-DD HH:MM:SS"

